# Fruta / hortaliza, semanal
# Inserts two new weekly price rows (573-574) into the "Frutilla" price
# history sheet, pushing the former rows 573:597 down to 575:599.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 573 (format/inherits from the row above,
# which keeps the date column's number format consistent with the rest of
# the table).
$ws.Rows.Item(573).Resize(2).EntireRow.Insert()

# New row 573
$ws.Range("A573").Value = 5
$ws.Range("B573").Value = "Macroferia Regional de Talca"
$ws.Range("C573").Value = "Maule"
$ws.Range("D573").Value = 44747
$ws.Range("E573").Value = 7
$ws.Range("F573").Value = "Fruta"
$ws.Range("G573").Value = 100101
$ws.Range("H573").Value = "Berries"
$ws.Range("I573").Value = 100112025
$ws.Range("J573").Value = "Frutilla"
$ws.Range("K573").Value = "Sin especificar"
$ws.Range("L573").Value = "Especial"
$ws.Range("M573").Value = 50
$ws.Range("N573").Value = 19000
$ws.Range("O573").Value = 19000
$ws.Range("P573").Value = 19000
$ws.Range("Q573").Value = "$/bandeja 7 kilos"
$ws.Range("R573").Value = "Provincia de Melipilla"
$ws.Range("S573").Value = 2714
$ws.Range("T573").Value = 7

# New row 574
$ws.Range("A574").Value = 5
$ws.Range("B574").Value = "Macroferia Regional de Talca"
$ws.Range("C574").Value = "Maule"
$ws.Range("D574").Value = 44747
$ws.Range("E574").Value = 7
$ws.Range("F574").Value = "Fruta"
$ws.Range("G574").Value = 100101
$ws.Range("H574").Value = "Berries"
$ws.Range("I574").Value = 100112025
$ws.Range("J574").Value = "Frutilla"
$ws.Range("K574").Value = "Sin especificar"
$ws.Range("L574").Value = "Segunda"
$ws.Range("M574").Value = 30
$ws.Range("N574").Value = 10000
$ws.Range("O574").Value = 10000
$ws.Range("P574").Value = 10000
$ws.Range("Q574").Value = "$/bandeja 7 kilos"
$ws.Range("R574").Value = "Provincia de Melipilla"
$ws.Range("S574").Value = 1429
$ws.Range("T574").Value = 7
